$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 581 (old rows 581-587 shift down to 585-591)
$ws.Rows("581:584").Insert()

# Row 581: new data (Patagonia, 1a nueva(o))
$ws.Cells.Item(581, 1).Value = 10
$ws.Cells.Item(581, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(581, 3).Value = "La Araucanía"
$ws.Cells.Item(581, 4).Value = 44595
$ws.Cells.Item(581, 5).Value = 9
$ws.Cells.Item(581, 6).Value = 100114001
$ws.Cells.Item(581, 7).Value = "Papa"
$ws.Cells.Item(581, 8).Value = "Patagonia"
$ws.Cells.Item(581, 9).Value = "1a nueva(o)"
$ws.Cells.Item(581, 10).Value = 500
$ws.Cells.Item(581, 11).Value = 7000
$ws.Cells.Item(581, 12).Value = 7000
$ws.Cells.Item(581, 13).Value = 7000
$ws.Cells.Item(581, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(581, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(581, 16).Value = 280
$ws.Cells.Item(581, 17).Value = 25
$ws.Cells.Item(581, 18).Value = "Hortaliza"

# Row 582: new data (Patagonia, 1a nueva(o))
$ws.Cells.Item(582, 1).Value = 10
$ws.Cells.Item(582, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(582, 3).Value = "La Araucanía"
$ws.Cells.Item(582, 4).Value = 44595
$ws.Cells.Item(582, 5).Value = 9
$ws.Cells.Item(582, 6).Value = 100114001
$ws.Cells.Item(582, 7).Value = "Papa"
$ws.Cells.Item(582, 8).Value = "Patagonia"
$ws.Cells.Item(582, 9).Value = "1a nueva(o)"
$ws.Cells.Item(582, 10).Value = 1000
$ws.Cells.Item(582, 11).Value = 6000
$ws.Cells.Item(582, 12).Value = 6000
$ws.Cells.Item(582, 13).Value = 6000
$ws.Cells.Item(582, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(582, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(582, 16).Value = 240
$ws.Cells.Item(582, 17).Value = 25
$ws.Cells.Item(582, 18).Value = "Hortaliza"

# Row 583: new data (Rodeo, 1a nueva(o))
$ws.Cells.Item(583, 1).Value = 10
$ws.Cells.Item(583, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(583, 3).Value = "La Araucanía"
$ws.Cells.Item(583, 4).Value = 44595
$ws.Cells.Item(583, 5).Value = 9
$ws.Cells.Item(583, 6).Value = 100114001
$ws.Cells.Item(583, 7).Value = "Papa"
$ws.Cells.Item(583, 8).Value = "Rodeo"
$ws.Cells.Item(583, 9).Value = "1a nueva(o)"
$ws.Cells.Item(583, 10).Value = 1200
$ws.Cells.Item(583, 11).Value = 7000
$ws.Cells.Item(583, 12).Value = 7000
$ws.Cells.Item(583, 13).Value = 7000
$ws.Cells.Item(583, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(583, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(583, 16).Value = 280
$ws.Cells.Item(583, 17).Value = 25
$ws.Cells.Item(583, 18).Value = "Hortaliza"

# Row 584: new data (Rodeo, 1a nueva(o))
$ws.Cells.Item(584, 1).Value = 10
$ws.Cells.Item(584, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(584, 3).Value = "La Araucanía"
$ws.Cells.Item(584, 4).Value = 44595
$ws.Cells.Item(584, 5).Value = 9
$ws.Cells.Item(584, 6).Value = 100114001
$ws.Cells.Item(584, 7).Value = "Papa"
$ws.Cells.Item(584, 8).Value = "Rodeo"
$ws.Cells.Item(584, 9).Value = "1a nueva(o)"
$ws.Cells.Item(584, 10).Value = 300
$ws.Cells.Item(584, 11).Value = 6000
$ws.Cells.Item(584, 12).Value = 6000
$ws.Cells.Item(584, 13).Value = 6000
$ws.Cells.Item(584, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(584, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(584, 16).Value = 240
$ws.Cells.Item(584, 17).Value = 25
$ws.Cells.Item(584, 18).Value = "Hortaliza"
